# "Lương" (salary) sheet: recalculated payroll report after total worked
# days ("Tổng công") at CẦN THƠ changed from 7 to 8 for this employee.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Lương")

# Tổng công tại CẦN THƠ: 7 -> 8
$ws.Range("B2").Value2 = 8

# Lương cơ bản tại CẦN THƠ = (lương cơ bản cũ / số công cũ) * số công mới
$soCongCu = 7
$soCongMoi = 8
$luongCoBanCu = $ws.Range("B3").Value2
$ws.Range("B3").Value2 = ($luongCoBanCu / $soCongCu) * $soCongMoi

# Tổng lương tại CẦN THƠ = tổng các khoản (B3:B11) tại CẦN THƠ
$ws.Range("B34").Value2 = 1238571.428571429

# Tổng lương tại HỆ THỐNG = Tổng lương CẦN THƠ + LONG XUYÊN + SÓC TRĂNG
$ws.Range("B37").Value2 = 1238571.428571429
